$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.902.82'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.03%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.635.66'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.14%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.76'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.07%  '
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.47%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '23.41'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +1.11%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.16%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0611'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -0.04%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.26%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.869.04'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.16%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.640.76'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.16%  '
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -0.71%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.563'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.81%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.40'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.52%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '27.919.84'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '228.98'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.15%  '
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +2.86%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0719'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.04%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.02%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.08%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.07'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -2.60%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.35%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '155.75'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +2.00%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.65%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.14%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.32%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +0.02%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.02%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.06%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +1.36%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +1.77%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.396.71'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -0.17%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.60'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +2.93%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.41%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.71%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.89%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.561'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.31%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.851'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -2.15%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.04%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.00'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -1.41%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.85'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +2.80%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '66.01'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -1.12%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.46'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.68%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.777.83'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.10%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -2.44%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '88.72'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +1.17%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +2.67%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.14%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.62'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.81%  '
